$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F8").Value = "68i"
$ws.Range("J5").Value = "oh!"
$ws.Range("J5").Select()
